$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Moved bass from John to Gordon" - rename the band member in B3
# (was "John", now "Gordon") while keeping everything else the same.
$ws.Range("B3").Value = "Gordon"

# Update the active selection to match the saved workbook state.
$ws.Range("B4").Select() | Out-Null
